$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their exact text representation
# (trailing zeros, percent signs, etc.) by formatting the target cells as Text
# before assigning the new values - mirrors how the source data was authored.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "328.94"
$ws.Range("E2").Value = "0.04%"
$ws.Range("D3").Value = "44.06"
$ws.Range("E3").Value = "-0.43%"
$ws.Range("D4").Value = "5.498"
$ws.Range("E4").Value = "-1.45%"
$ws.Range("D5").Value = "0.08081"
$ws.Range("E5").Value = "0.26%"
$ws.Range("D6").Value = "2.022"
$ws.Range("E6").Value = "5.27%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.411"
$ws.Range("E7").Value = "2.94%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "0.9508"
$ws.Range("E8").Value = "0.20%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1126"
$ws.Range("E9").Value = "-7.30%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1881"
$ws.Range("E10").Value = "2.09%"
$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").Value = "10.12"
$ws.Range("E11").Value = "1.17%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.1001"
$ws.Range("E12").Value = "3.29%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.04828"
$ws.Range("E13").Value = "10.58%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.1058"
$ws.Range("E14").Value = "-0.78%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001260"
$ws.Range("E15").Value = "-1.29%"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "0.04080"
$ws.Range("E16").Value = "-3.20%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.005988"
$ws.Range("E17").Value = "0.65%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.371"
$ws.Range("E18").Value = "-0.68%"
$ws.Range("D19").Value = "2.602"
$ws.Range("E19").Value = "1.08%"
$ws.Range("E20").Value = "-4.89%"
$ws.Range("D21").Value = "0.1400"
$ws.Range("E21").Value = "-1.89%"
$ws.Range("E22").Value = "2.74%"
$ws.Range("D23").Value = "0.001308"
$ws.Range("E23").Value = "5.22%"
$ws.Range("D24").Value = "0.004363"
$ws.Range("E24").Value = "0.66%"
$ws.Range("D25").Value = "0.0001249"
$ws.Range("E25").Value = "5.14%"
$ws.Range("D26").Value = "0.0003740"
$ws.Range("E26").Value = "-6.01%"
$ws.Range("D38").Value = "0.02616"
$ws.Range("E38").Value = "-2.12%"
$ws.Range("D39").Value = "0.05655"
$ws.Range("E39").Value = "2.33%"
$ws.Range("D40").Value = "0.007589"
$ws.Range("E40").Value = "0.55%"
$ws.Range("D41").Value = "0.1403"
$ws.Range("E41").Value = "-0.16%"
$ws.Range("D42").Value = "0.007344"
$ws.Range("E42").Value = "-11.79%"
$ws.Range("D43").Value = "0.001984"
$ws.Range("E43").Value = "-1.42%"
$ws.Range("D44").Value = "0.008257"
$ws.Range("E44").Value = "-7.10%"
$ws.Range("D45").Value = "0.00007079"
$ws.Range("E45").Value = "-0.36%"
$ws.Range("E46").Value = "0.12%"
$ws.Range("D47").Value = "0.0005800"
$ws.Range("E47").Value = "-0.20%"
$ws.Range("D48").Value = "0.003498"
$ws.Range("E48").Value = "53.98%"
$ws.Range("D49").Value = "0.003508"
$ws.Range("E49").Value = "24.00%"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "0.12%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "0.12%"

Write-Output "applied crypto price update"
